# Add a new "2020" column (K) to the municipal-solid-waste-disposal table,
# mirroring the existing 2013-2019 columns (A:J) that are already on the sheet.
# This replicates the common "copy the last year's column, paste it next to
# itself, then type in the new year's figures" editing workflow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy column J (2019, rows 2-8 incl. the header/border rows) into column K
#    so the new column inherits the same fonts/borders/number-formats as the
#    rest of the table.
$ws.Range("J2:J8").Copy()
$ws.Range("K2:K8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2) Overwrite the pasted values with the real 2020 figures.
$ws.Range("K3").Value2 = 2020
$ws.Range("K4").Value2 = 0
$ws.Range("K5").Value2 = 48.2
$ws.Range("K6").Value2 = 19.3
$ws.Range("K7").Value2 = 24.2
$ws.Range("K8").Value2 = 8.3000000000000007

# 3) K5 (the "Мусоропровод" / garbage-chute row for 2020) was additionally
#    reformatted by hand: one decimal place, right aligned.
$ws.Range("K5").NumberFormat = "0.0"
$ws.Range("K5").HorizontalAlignment = -4152
$ws.Range("K5").Font.Name = "Times New Roman"
$ws.Range("K5").Font.Size = 9

# 4) The header row got shorter once the sheet was touched again.
$ws.Rows.Item(1).RowHeight = 63.75

# 5) Leave the selection where the author left it after editing.
$ws.Range("J22").Select()
